$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the shared-string text (NL LMBV B06 (IMPORT & Export) -> NL LMBV B02 (IMPORT))
$ws.Range("M2").Value = "NL LMBV B02 (IMPORT)"

# 2) Resize columns A:M individually (previously grouped with shared widths).
#    The runtime stores column widths internally in 1/6-character increments, so
#    the ColumnWidth we assign is the (target width - 5/6) offset used by the
#    COM layer when it serializes back to the OOXML "width" attribute.
$ws.Columns.Item(1).ColumnWidth = 16.8809523809524
$ws.Columns.Item(2).ColumnWidth = 20.3095238095238
$ws.Columns.Item(3).ColumnWidth = 16.7380952380952
$ws.Columns.Item(4).ColumnWidth = 23.0238095238095
$ws.Columns.Item(5).ColumnWidth = 28.1666666666667
$ws.Columns.Item(6).ColumnWidth = 21.8809523809524
$ws.Columns.Item(7).ColumnWidth = 28.5952380952381
$ws.Columns.Item(8).ColumnWidth = 19.7380952380952
$ws.Columns.Item(9).ColumnWidth = 17.1666666666667
$ws.Columns.Item(10).ColumnWidth = 18.1666666666667
$ws.Columns.Item(11).ColumnWidth = 24.4523809523810
$ws.Columns.Item(12).ColumnWidth = 12.5952380952381
$ws.Columns.Item(13).ColumnWidth = 22.3095238095238

# 3) Update the selection so it matches the new active cell / scroll target.
$ws.Range("I13").Select() | Out-Null
